# Automatic map update (2025-08-11 07:12:43)
#
# The source "NEW" sheet lost its earliest pending ticket (Caso 5825,
# "PAZ, GRAL. AV. 5602", row 36): it was resolved/removed upstream, so the
# whole row is deleted and every row below it shifts up by one. The sheet
# therefore shrinks from A1:P61 to A1:P60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(36).Delete()
